$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "Compact explicit multi-path routing for LEO satellite networks"
$ws.Range("B5").Value = "HPSR"
$ws.Range("C5").Value = "1.对路径进行编码，用全局的编码替代header里面的路径，以此减少传输的消耗。`n2.如果路径编码无效，则使用默认最短路径"
$ws.Range("D5").Value = "1.指出ISL由于天线转向等因素在纬度70度左右会断开连接，"
$ws.Range("E5").Value = "1.过于简单，只是减少了header的大小，减少的传输消耗的效果有待商榷"

# Row 6 (publication entered before title, matching original authoring order)
$ws.Range("B6").Value = "Computer Networks"
$ws.Range("A6").Value = "Distributed on-demand routing for LEO satellite systems"
$ws.Range("C6").Value = "1.提出收缩广播的区域，在小范围内通过广播寻找最短路径。`n2.寻路方式：通过广播探测包寻找最短路径，每条路径有一定的有效期"
$ws.Range("D6").Value = "1.同样指出ISL在极低附近会断开与相邻轨道的连接`n2.缩小广播范围确实能减少寻路开销"
$ws.Range("E6").Value = "1.广播域限缩在源/目的节点内的矩形内，当长距离传输时，广播域会非常大，`n2.考虑到排队时延，局部网络拥塞可能导致广播域内找不到时延低的路径`n3."

# Apply styles to match rows 3/4 pattern (B col centered, C/D/E wrap text)
$ws.Range("B4").Copy()
$ws.Range("B5:B6").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Range("C4:E4").Copy()
$ws.Range("C6:E6").PasteSpecial(-4122)

# Row heights
$ws.Rows.Item(5).RowHeight = 100
$ws.Rows.Item(6).RowHeight = 92

# Update view/selection to match new dimension
$ws.Range("E6").Select()
